$wb = $excel.ActiveWorkbook

# Sheets: 1=Bus, 2=Gen, 3=Branch
$busSheet = $wb.Worksheets.Item(1)
$genSheet = $wb.Worksheets.Item(2)

# Gen sheet ("B matrix") gets a new row of data (row 10), extending the
# used range from A1:U9 to A1:U10.
$genSheet.Cells.Item(10, 1).Value = 9
for ($col = 2; $col -le 21; $col++) {
    $genSheet.Cells.Item(10, $col).Value = 0
}

# Update the selection on the Gen sheet to reflect the newly added rows,
# without making Gen the active/selected tab.
$null = $genSheet.Range("B9:U10").Select()

# Make the Bus sheet the selected/active tab (previously Gen was active).
$null = $busSheet.Activate()
